$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.344.50"
$ws.Range("E2").Value = "  -1.17%  "

$ws.Range("D3").Value = "2.371.30"
$ws.Range("E3").Value = "  +5.56%  "

$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("E5").Value = "  +1.23%  "

$ws.Range("D6").Value = "'232.55"
$ws.Range("E6").Value = "  +0.97%  "

$ws.Range("D7").Value = "'68.63"
$ws.Range("E7").Value = "  +6.65%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("E9").Value = "  +2.37%  "

$ws.Range("D10").Value = "'0.0952"
$ws.Range("E10").Value = "  -2.47%  "

$ws.Range("D11").Value = "'56.90"
$ws.Range("E11").Value = "  -0.18%  "

$ws.Range("E12").Value = "  -0.74%  "

$ws.Range("D13").Value = "2.726.18"
$ws.Range("E13").Value = "  +5.60%  "

$ws.Range("E14").Value = "  -0.54%  "

$ws.Range("D15").Value = "'15.61"
$ws.Range("E15").Value = "  +0.44%  "

$ws.Range("E16").Value = "  +3.57%  "

$ws.Range("D17").Value = "'0.846"
$ws.Range("E17").Value = "  +1.87%  "

$ws.Range("D18").Value = "2.379.63"
$ws.Range("E18").Value = "  +5.47%  "

$ws.Range("D19").Value = "43.374.07"
$ws.Range("E19").Value = "  -0.70%  "

$ws.Range("D20").Value = "0.0₃0984"
$ws.Range("E20").Value = "  -0.14%  "

$ws.Range("D21").Value = "'6.28"
$ws.Range("E21").Value = "  +4.21%  "

$ws.Range("D22").Value = "'73.92"
$ws.Range("E22").Value = "  +1.09%  "

$ws.Range("D23").Value = "'249.35"
$ws.Range("E23").Value = "  -0.42%  "

$ws.Range("D24").Value = "'3.91"
$ws.Range("E24").Value = "  +17.46%  "

$ws.Range("E25").Value = "  -0.01%  "

$ws.Range("E26").Value = "  +1.27%  "

$ws.Range("E27").Value = "  -1.35%  "

$ws.Range("D28").Value = "'10.01"
$ws.Range("E28").Value = "  +0.35%  "

$ws.Range("D29").Value = "'22.44"
$ws.Range("E29").Value = "  +7.94%  "

$ws.Range("D30").Value = "'175.30"
$ws.Range("E30").Value = "  +2.76%  "

$ws.Range("E31").Value = "  +9.12%  "

$ws.Range("E32").Value = "  -5.93%  "

$ws.Range("E33").Value = "  +0.75%  "

$ws.Range("E34").Value = "  +5.04%  "

$ws.Range("D35").Value = "'0.0695"
$ws.Range("E35").Value = "  -0.99%  "

$ws.Range("D36").Value = "'5.03"
$ws.Range("E36").Value = "  +3.46%  "

$ws.Range("E37").Value = "  +11.03%  "

$ws.Range("E38").Value = "  +0.96%  "

$ws.Range("E39").Value = "  -1.62%  "

$ws.Range("E40").Value = "  -2.29%  "

$ws.Range("E41").Value = "  +10.52%  "

$ws.Range("E42").Value = "  -0.10%  "

$ws.Range("D43").Value = "'17.97"
$ws.Range("E43").Value = "  +4.05%  "

$ws.Range("E44").Value = "  +9.12%  "

$ws.Range("D45").Value = "'99.42"
$ws.Range("E45").Value = "  +2.23%  "

$ws.Range("E46").Value = "  +1.28%  "

$ws.Range("D47").Value = "'0.0954"
$ws.Range("E47").Value = "  -0.87%  "

$ws.Range("E48").Value = "  -0.42%  "

$ws.Range("D49").Value = "1.448.13"
$ws.Range("E49").Value = "  +1.08%  "

$ws.Range("B50").Value = "TerraClassic"
$ws.Range("C50").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D50").Value = "'0.000206"
$ws.Range("E50").Value = "  -8.26%  "

$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.598.55"
$ws.Range("E51").Value = "  +5.82%  "
